$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("posters")

# --- Update existing citation cells (append poster-board / session numbers) ---
$ws.Range("D2").Value = 'Ziman K, Lee MR, Martinez AR, Manning JR (2019) Volitional Attention Modulates Memory Encoding and Retrieval. <em>Society for Neuroscience.</em> Chicago, IL. 792.22/LLL18.'
$ws.Range("D3").Value = 'Fitzpatrick PC, Heusser AC, Manning JR (2019) Exploring the evolving geometric structure of experiences and memories. <em>Society for Neuroscience</em>. Chicago, IL. 423.16/BB14.'
$ws.Range("D4").Value = 'Owen LLW, Manning JR (2019) Understanding complexity and interactivity of brain patterns in naturalistic processing. <em>Society for Neuroscience.</em> Chicago, IL. 248.17/Z41.'
$ws.Range("D6").Value = 'Ziman K, Lee MR, Martinez AR, Manning JR (2018) Volitional Attention Modulates Memory Encoding and Retrieval. <em>Society for Neuroscience.</em> San Diego, CA. 792.22/LLL18.'
$ws.Range("D7").Value = 'Fitzpatrick PC, Heusser AC, Manning JR (2018) Mapping between naturalistic experience and verbal recall. <em>Society for Neuroscience.</em> San Diego, CA. 086.10/HHH33.'
$ws.Range("D12").Value = 'Manning JR, Ziman K, Heusser AC (2017) Efficient Learning: Manipulating context to enhance (or diminish) memory. <em>Society for Neuroscience.</em> Washington, DC. 339.08/UU42.'
$ws.Range("D13").Value = 'Heusser AC, Ziman K, Owen LLW, Manning JR (2017) HyperTools: A python toolbox for gaining geometric insights into high-dimensional data. <em>Society for Neuroscience.</em> Washington, DC. 721.22/WW28.'
$ws.Range("D14").Value = 'Owen LLW, Manning JR (2017) A Gaussian process model of human ECoG data. <em>Society for Neuroscience.</em> Washington, DC. 093.04/UU78.'
$ws.Range("D15").Value = 'Ziman K, Heusser AC, Manning JR (2017) Effects of study context on recall organization. <em>Society for Neuroscience.</em> Washington, DC. 803.07/UU14.'

# --- Add PDF links to existing Wetterhahn 2023/2024 rows ---
$ws.Range("E44").Value = '[<a href="data/pdfs/CarsEtal24.pdf" target="_blank">pdf</a>]'
$ws.Range("E45").Value = '[<a href="data/pdfs/JhaEtal23.pdf" target="_blank">pdf</a>]'

# --- New rows 46-50 (title + citation only) ---
$ws.Range("B46").Value = 'Episodic memory: Mental time travel or a quantum "memory wave" function?'
$ws.Range("D46").Value = 'Manning JR (2019) Episodic memory: Mental time travel or a quantum "memory wave" function? <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA.'
$ws.Range("B47").Value = 'A geometric approach to modeling knowledge and learning from Khan Academy course videos'
$ws.Range("D47").Value = 'Fitzpatrick PC, Heusser AC, Manning JR (2022) A geometric approach to modeling knowledge and learning from Khan Academy course videos. <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA.'
$ws.Range("B48").Value = 'Temporal asymmetries in narrative events'
$ws.Range("D48").Value = 'Xu X (2022) Temporal asymmetries in narrative events. <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA.'
$ws.Range("B49").Value = 'Modeling the knowledge asymmetry of the past and the future'
$ws.Range("D49").Value = 'Xu X (2024) Modeling the knowledge asymmetry of the past and the future. <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA.'
$ws.Range("B50").Value = 'Shared representational geometry as an explanation for cross-subject prediction of place cell data from the rodent hippocampus'
$ws.Range("D50").Value = 'Chen HT, Manning JR, van der Meer MAA (2019) Shared representational geometry as an explanation for cross-subject prediction of place cell data from the rodent hippocampus. <em>Society for Neuroscience</em>, Chicago, IL. 084.04/Y28.'

# --- New rows 51-59 (title + citation + pdf link) ---
$ws.Range("B51").Value = 'Translating neurophysiological recordings into dynamic estimates of conceptual knowledge and learning'
$ws.Range("D51").Value = 'Peng K, Carstensen D, Parigela S, Shah O, Wingo A, Liu A, Maina J, Dampal K, Manning JR (2025) Translating neurophysiological recordings into dynamic estimates of conceptual knowledge and learning. <em>Karen E. Wetterhahn Science Symposium</em>, Hanover, NH.'
$ws.Range("E51").Value = '[<a href="data/pdfs/PengEtal25.pdf" target="_blank">pdf</a>]'
$ws.Range("B52").Value = 'Higher-order interactions between brain regions are better at profiling tasks'
$ws.Range("D52").Value = 'Saggar M, Betzel R, Manning JR, Liegeois R, Sporns O, Petri G (2022) Higher-order interactions between brain regions are better at profiling tasks. <em>Organization for Human Brain Mapping</em>, Glasgow, Scotland.'
$ws.Range("E52").Value = '[<a href="data/pdfs/SaggEtal22.pdf" target="_blank">pdf</a>]'
$ws.Range("B53").Value = 'Cognitive markers of mental health'
$ws.Range("D53").Value = 'Jain S, Schreder N, Fitzpatrick PC, Ziman K, Manning JR (2022) Cognitive markers of mental health. <em>Conference on Cognitive Computational Neuroscience</em>, San Francisco, CA.'
$ws.Range("E53").Value = '[<a href="data/pdfs/JainEtal22.pdf" target="_blank">pdf</a>]'
$ws.Range("B54").Value = 'A Gaussian process model of human ECoG data'
$ws.Range("D54").Value = 'Owen LLW, Muntianu TA, Heusser AC, Manning JR (2020) A Gaussian process model of human ECoG data. <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA (virtual).'
$ws.Range("E54").Value = '[<a href="data/pdfs/OwenEtal20.pdf" target="_blank">pdf</a>]'
$ws.Range("B55").Value = 'Understanding brain pattern complexity and interactivity in naturalistic processing'
$ws.Range("D55").Value = 'Owen LLW, Manning JR (2020) Understanding brain pattern complexity and interactivity in naturalistic processing. <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA (virtual).'
$ws.Range("E55").Value = '[<a href="data/pdfs/OwenMann20.pdf" target="_blank">pdf</a>]'
$ws.Range("B56").Value = 'Extrapolating the unobserved past and future in other people''s autobiographical timelines'
$ws.Range("D56").Value = 'Xu X, Zhu Z, Manning JR (2021) Extrapolating the unobserved past and future in other people''s autobiographical timelines. <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA.'
$ws.Range("E56").Value = '[<a href="data/pdfs/XuEtal21.pdf" target="_blank">pdf</a>]'
$ws.Range("B57").Value = 'Why we know more about the past: insights from statistical mechanics'
$ws.Range("D57").Value = 'Xu X, Manning JR (2023) Why we know more about the past: insights from statistical mechanics. <em>Context and Episodic Memory Symposium</em>, Orlando, FL.'
$ws.Range("E57").Value = '[<a href="data/pdfs/XuMann23.pdf" target="_blank">pdf</a>]'
$ws.Range("B58").Value = 'Temporal asymmetries in cued recall of naturalistic events'
$ws.Range("D58").Value = 'Xu X, Manning JR (2025) Temporal asymmetries in cued recall of naturalistic events. <em>Context and Episodic Memory Symposium</em>, Philadelphia, PA.'
$ws.Range("E58").Value = '[<a href="data/pdfs/XuMann25.pdf" target="_blank">pdf</a>]'
$ws.Range("B59").Value = 'Unexpected false feelings of familiarity about faces are associated with increased pupil dilations'
$ws.Range("D59").Value = 'Ziman K, Manning JR (2021) Unexpected false feelings of familiarity about faces are associated with increased pupil dilations. <em>Society for Neuroscience</em>, Virtual.'
$ws.Range("E59").Value = '[<a href="data/pdfs/ZimaMann21.pdf" target="_blank">pdf</a>]'
